$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('vida rural, vida no campo', 3.378378378378379, 6.756756756756757, 0.6666666666666667, 0),
    @('saúde', 2.702702702702703, 7.432432432432432, 1.333333333333333, 0),
    @('sistema penitenciário', 0.6756756756756757, 0.6756756756756757, 0.6666666666666667, 0),
    @('recreação, lazer e entretenimento', 24.32432432432433, 54.05405405405406, 0.6666666666666667, 0),
    @('outros', 27.7027027027027, 68.24324324324324, 0, 0),
    @('nacionalidade e nacionalismo', 4.054054054054054, 7.432432432432432, 0, 0),
    @('raça e racismo', 0.6756756756756757, 3.378378378378379, 1.333333333333333, 0),
    @('maternidade', 1.351351351351351, 2.027027027027027, 0, 0),
    @('linguagem e representação', 12.16216216216216, 36.48648648648648, 0, 0),
    @('questão agrária e territorial', 4.72972972972973, 11.48648648648649, 4.666666666666667, 0.6666666666666667),
    @('sociedade e consumo', 10.81081081081081, 20.27027027027027, 1.333333333333333, 0.6666666666666667),
    @('pessoas com deficiência', 1.351351351351351, 0.6756756756756757, 0.6666666666666667, 0.6666666666666667),
    @('violência', 4.054054054054054, 12.83783783783784, 1.333333333333333, 0.6666666666666667),
    @('sonho e fantasia', 2.027027027027027, 2.027027027027027, 0.6666666666666667, 0.6666666666666667),
    @('moradia e habitação', 0.6756756756756757, 3.378378378378379, 4, 0.6666666666666667),
    @('violências e preconceitos de gênero', 2.702702702702703, 5.405405405405405, 1.333333333333333, 0.6666666666666667),
    @('amizade', 0.6756756756756757, 1.351351351351351, 0, 0.6666666666666667),
    @('família', 5.405405405405405, 12.16216216216216, 6, 0.6666666666666667),
    @('distopia, ficção científica e futuros imaginados', 11.48648648648649, 10.81081081081081, 1.333333333333333, 1.333333333333333),
    @('desinformação, populismo e polarização', 2.702702702702703, 6.081081081081082, 1.333333333333333, 1.333333333333333),
    @('territorialidade e colonialismo', 5.405405405405405, 10.81081081081081, 3.333333333333333, 1.333333333333333),
    @('tecnologia, inovação e sociedade', 6.756756756756757, 10.81081081081081, 1.333333333333333, 1.333333333333333),
    @('sonoridade e paisagem sonora', 8.783783783783784, 13.51351351351351, 2.666666666666667, 1.333333333333333),
    @('mobilidade', 4.72972972972973, 10.13513513513514, 5.333333333333334, 2),
    @('gênero e sexualidade', 12.16216216216216, 17.56756756756757, 2, 2.666666666666667),
    @('direitos humanos', 4.054054054054054, 12.83783783783784, 0, 3.333333333333333),
    @('saudade, luto ou perda', 7.432432432432432, 13.51351351351351, 0.6666666666666667, 3.333333333333333),
    @('religião, espiritualidade e cosmologias', 4.72972972972973, 10.81081081081081, 2, 4),
    @('trabalho e ofício', 10.13513513513514, 20.27027027027027, 5.333333333333334, 4),
    @('política', 12.16216216216216, 31.08108108108108, 6, 4),
    @('solidão', 1.351351351351351, 4.054054054054054, 2, 4),
    @('biografia', 2.027027027027027, 6.756756756756757, 0, 4),
    @('alimentação e tratamentos tradicionais', 4.72972972972973, 4.72972972972973, 2, 4),
    @('corpo, performance e expressão', 4.054054054054054, 12.83783783783784, 1.333333333333333, 4),
    @('dinâmica urbana', 24.32432432432433, 42.56756756756756, 9.333333333333334, 4),
    @('estruturas sociais e econômicas', 26.35135135135135, 52.70270270270269, 4.666666666666667, 4),
    @('crises e desastres ambientais e sociais', 4.054054054054054, 6.756756756756757, 3.333333333333333, 4.666666666666667),
    @('cultural', 34.45945945945946, 56.75675675675676, 7.333333333333333, 4.666666666666667),
    @('vida cotidiana', 5.405405405405405, 14.86486486486486, 2, 5.333333333333334),
    @('memória e patrimônio', 16.89189189189189, 39.18918918918919, 6, 5.333333333333334),
    @('educação e socialização', 14.18918918918919, 33.10810810810811, 6.666666666666667, 6),
    @('resistência e luta', 6.756756756756757, 12.83783783783784, 4.666666666666667, 6),
    @('vida afetiva', 49.32432432432432, 69.5945945945946, 11.33333333333333, 6),
    @('pandemia', 16.21621621621622, 20.94594594594595, 14.66666666666667, 7.333333333333333),
    @('saúde mental', 14.86486486486486, 29.72972972972973, 8, 8),
    @('mulher, feminino e feminismo', 4.054054054054054, 6.081081081081082, 6.666666666666667, 8),
    @('reflexão', 2.702702702702703, 10.13513513513514, 2, 8),
    @('povos originários e comunidades tradicionais', 13.51351351351351, 16.21621621621622, 3.333333333333333, 8),
    @('poesia e ensaio', 37.83783783783784, 64.86486486486487, 12, 11.33333333333333),
    @('ambiental', 16.21621621621622, 22.97297297297298, 11.33333333333333, 12),
    @('arte', 3.378378378378379, 6.756756756756757, 3.333333333333333, 13.33333333333333)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}

Write-Host "Updated rows 2-52"